$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 05:46"

# Update country names (column A) to reflect the re-sorted ranking
$ws.Range("A123").Value = "Honduras"
$ws.Range("A124").Value = "Mayotte"
$ws.Range("A128").Value = "Polinesia Francesa"
$ws.Range("A129").Value = "Kenia"
$ws.Range("A137").Value = "Madagascar"
$ws.Range("A138").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("A144").Value = "Mongolia"
$ws.Range("A145").Value = "Nueva Caledonia"
$ws.Range("A146").Value = "El Salvador"
$ws.Range("A148").Value = "Uganda"
$ws.Range("A149").Value = "San Martin (Parte Francesa)"
$ws.Range("A150").Value = "Dominica"
$ws.Range("A151").Value = "Surinam"
$ws.Range("A153").Value = "Seychelles"
$ws.Range("A154").Value = "Namibia"
$ws.Range("A155").Value = "Bermudas"
$ws.Range("A156").Value = "Benin"
$ws.Range("A157").Value = "Islas Caimanes"
$ws.Range("A158").Value = "Curazao"
$ws.Range("A159").Value = "Gabon"
$ws.Range("A160").Value = "Fiyi"
$ws.Range("A162").Value = "Guyana"
$ws.Range("A163").Value = "Groenlandia"
$ws.Range("A165").Value = "Congo"
$ws.Range("A166").Value = "Santa Sede"
$ws.Range("A168").Value = "San Bartolome"
$ws.Range("A171").Value = "Republica del Chad"
$ws.Range("A172").Value = "Republica de Yibuti"
$ws.Range("A173").Value = "Republica de Africa Central"
$ws.Range("A174").Value = "Zambia"
$ws.Range("A175").Value = "Santa Lucia"
$ws.Range("A176").Value = "Birmania"
$ws.Range("A177").Value = "Niger"
$ws.Range("A180").Value = "Cabo Verde"
$ws.Range("A181").Value = "Gambia"
$ws.Range("A183").Value = "Zimbabue"
$ws.Range("A185").Value = "Mauritania"
$ws.Range("A186").Value = "Laos"
$ws.Range("A187").Value = "Nicaragua"
$ws.Range("A188").Value = "San Martin (Parte Holandesa)"
$ws.Range("A190").Value = "Somalia"
$ws.Range("A191").Value = "Siria"
$ws.Range("A192").Value = "San Vicente y las Granadinas"
$ws.Range("A194").Value = "Libia"
$ws.Range("A195").Value = "Papua Nueva Guinea"
$ws.Range("A196").Value = "Timor Oriental"
$ws.Range("A197").Value = "Eritrea"
$ws.Range("A198").Value = "Montserrat"
$ws.Range("A199").Value = "Granada"
$ws.Range("A200").Value = "Islas Turcas y Caicos"

# Update updated case-count figures for the affected rows
$ws.Range("C123").Value = 6
$ws.Range("C124").Value = 0
$ws.Range("C146").Value = 4
$ws.Range("B148").Value = 9
$ws.Range("E148").Value = 9
$ws.Range("B149").Value = 8
$ws.Range("E149").Value = 8
$ws.Range("D153").Value = 0
$ws.Range("E153").Value = 7
$ws.Range("B154").Value = 7
$ws.Range("D154").Value = 2
$ws.Range("E154").Value = 5
$ws.Range("E156").Value = 6
$ws.Range("H156").Value = 0
$ws.Range("B159").Value = 6
$ws.Range("H159").Value = 1
$ws.Range("C160").Value = 1
$ws.Range("E160").Value = 5
$ws.Range("H160").Value = 0
$ws.Range("D162").Value = 0
$ws.Range("E162").Value = 4
$ws.Range("H162").Value = 1
$ws.Range("B163").Value = 5
$ws.Range("D163").Value = 2
$ws.Range("E163").Value = 3
